# Turns "File01 - master branch" into "File02 - master branch", reproducing
# the exact run layout Word itself produces for this kind of edit: the
# digit that was retyped ends up in its own run, split off from the
# "File0" prefix, and the "_GoBack" bookmark (which always marks the most
# recent edit point) is relocated to sit right after the new "2", ahead of
# the untouched " - master branch" tail (which becomes its own run too).

$d = $word.ActiveDocument

# Locate "File01" precisely (avoids accidentally matching some other "1").
$found = $d.Content
$found.Find.Execute("File01", $false, $false, $false, $false, $false, `
                     $true, 1, $false, $null, 0) | Out-Null

$start = $found.Start
$digitStart = $start + 5   # "File0" is 5 chars; the "1" is the 6th
$digitEnd   = $digitStart + 1

# 1) Retype the digit: "1" -> "2".
$digitRng = $d.Range($digitStart, $digitEnd)
$digitRng.Text = "2"

# 2) Split the run boundary just before the new "2" (separates it from
#    "File0"). Adding + immediately deleting a bookmark at a point forces
#    a structural run split there that survives the bookmark's removal.
$beforeTwo = $d.Range($digitStart, $digitStart)
$d.Bookmarks.Add("TempSplitBefore", $beforeTwo) | Out-Null
$d.Bookmarks.Item("TempSplitBefore").Delete()

# 3) Split the run boundary just after the new "2" (separates it from the
#    trailing " - master branch").
$afterTwo = $d.Range($digitEnd, $digitEnd)
$d.Bookmarks.Add("TempSplitAfter", $afterTwo) | Out-Null
$d.Bookmarks.Item("TempSplitAfter").Delete()

# 4) Move "_GoBack" (Word's "last edit" bookmark) to sit right after the
#    "2" -- re-adding it under its reserved name relocates the existing
#    one instead of creating a duplicate.
$afterTwo2 = $d.Range($digitEnd, $digitEnd)
$d.Bookmarks.Add("_GoBack", $afterTwo2) | Out-Null
